# CIV-8879 Removing correspondence address for the 2nd defendant
#
# The 2nd defendant's details table contains a "Correspondence address"
# row (label cell + a large block of <<cs_{...}>> merge fields for the
# representative's service address) that the 1st defendant's equivalent
# table never had. This row must be removed so both tables share the
# same shape again (Name / Address / Date of birth / Telephone / Email).

$d = $word.ActiveDocument

$targetRow = $null

foreach ($t in $d.Tables) {
    foreach ($r in $t.Rows) {
        $labelCell = $r.Cells.Item(1)
        $labelText = $labelCell.Range.Text.Trim().TrimEnd([char]7, [char]13)

        if ($labelText -eq "Correspondence address") {
            $targetRow = $r
            break
        }
    }
    if ($targetRow -ne $null) {
        break
    }
}

if ($targetRow -ne $null) {
    $targetRow.Delete()
} else {
    throw "Could not find the 'Correspondence address' row to delete."
}
